# The NATMI pipeline was re-run with updated TPM input data. The "Sending
# cluster" for this Il17b-Il17rb edge table changes from "MuSCs" to "ECs"
# for every data row (A2:A5), and all of the downstream per-row statistics
# that depend on the sending-cluster's expression (ligand detection rate,
# average/total expression, derived specificities, edge weights, ...) are
# refreshed with the newly computed numbers. The "Target cluster" column
# (D2:D5) keeps the same per-row values (ECs, FAPs, MuSCs, Resolving-Mac).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: Sending cluster -------------------------------------------
$ws.Range("A2:A5").Value = "ECs"

# --- Column D: Target cluster (unchanged text, kept explicit) ------------
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("D5").Value = "Resolving-Mac"

# --- Columns E-H: same new values across all four rows -------------------
$ws.Range("E2:E5").Value = 1
$ws.Range("F2:F5").Value = 0.3333333333333333
$ws.Range("G2:G5").Value = 0.128774
$ws.Range("H2:H5").Value = 0.386322

# --- Row 2 (Target cluster = ECs) -----------------------------------------
$ws.Range("M2").Value = 2.060782666666667
$ws.Range("N2").Value = 6.182348
$ws.Range("O2").Value = 0.2355590467032963
$ws.Range("P2").Value = 0.2355590467032963
$ws.Range("Q2").Value = 0.2653752271173334
$ws.Range("R2").Value = 2.388377044056
$ws.Range("S2").Value = 0.2355590467032963
$ws.Range("T2").Value = 0.2355590467032963

# --- Row 3 (Target cluster = FAPs) ----------------------------------------
$ws.Range("O3").Value = 0.3579488968516726
$ws.Range("P3").Value = 0.3579488968516725
$ws.Range("Q3").Value = 0.4032567253426666
$ws.Range("R3").Value = 3.629310528083999
$ws.Range("S3").Value = 0.3579488968516726
$ws.Range("T3").Value = 0.3579488968516725

# --- Row 4 (Target cluster = MuSCs) ---------------------------------------
$ws.Range("M4").Value = 2.215918666666667
$ws.Range("N4").Value = 6.647756
$ws.Range("O4").Value = 0.2532919638422357
$ws.Range("P4").Value = 0.2532919638422357
$ws.Range("Q4").Value = 0.2853527103813333
$ws.Range("R4").Value = 2.568174393432
$ws.Range("S4").Value = 0.2532919638422357
$ws.Range("T4").Value = 0.2532919638422357

# --- Row 5 (Target cluster = Resolving-Mac) -------------------------------
$ws.Range("M5").Value = 1.340267333333333
$ws.Range("N5").Value = 4.020802
$ws.Range("O5").Value = 0.1532000926027954
$ws.Range("P5").Value = 0.1532000926027954
$ws.Range("Q5").Value = 0.1725915855826667
$ws.Range("R5").Value = 1.553324270244
$ws.Range("S5").Value = 0.1532000926027954
$ws.Range("T5").Value = 0.1532000926027954
